# Auto-generated script applying the 'Add generated test cases from Claude' edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

# Row 6 (TC001)
$ws.Range("B6").Value = "TC001"
$ws.Range("C6").Value = "MFP is powered on and in normal boot-up mode"
$ws.Range("D6").Value = "Verify SCC function can be enabled/disabled through Self-diagnostic boot-up mode"
$ws.Range("E6").Value = "1. Enter Self-diagnostic boot-up mode`n2. Set 08-3820 to 1 (Enable)`n3. Reboot the MFP`n4. Verify SCC function is enabled`n5. Enter Self-diagnostic boot-up mode again`n6. Set 08-3820 to 0 (Disable)`n7. Reboot the MFP"
$ws.Range("F6").Value = "SCC function should be enabled when 08-3820 is set to 1 and disabled when set to 0"
$ws.Range("H6").Value = "Default value should be `"Disable`""

# Row 7 (TC002)
$ws.Range("B7").Value = "TC002"
$ws.Range("C7").Value = "MFP is powered on and in normal boot-up mode"
$ws.Range("D7").Value = "Verify SCC function cannot be enabled/disabled from TA or STAGE2 IF"
$ws.Range("E7").Value = "1. Log in to TopAccess as administrator`n2. Check if there is any option to enable/disable SCC`n3. Try to enable/disable SCC via STAGE2 IF"
$ws.Range("F7").Value = "No option to enable/disable SCC should be available in TopAccess or via STAGE2 IF"
$ws.Range("H7").Value = ""

# Row 8 (TC003)
$ws.Range("B8").Value = "TC003"
$ws.Range("C8").Value = "MFP is powered on and in normal boot-up mode"
$ws.Range("D8").Value = "Verify SCC Proxy Server Setting can be configured"
$ws.Range("E8").Value = "1. Enter Self-diagnostic boot-up mode`n2. Configure SCC Proxy Server settings (08-3822 to 08-3826)`n3. Reboot the MFP`n4. Verify proxy settings are applied"
$ws.Range("F8").Value = "Proxy server settings should be configurable only through Self-diagnostic mode"
$ws.Range("H8").Value = "Default value should be `"Disable`""

# Row 9 (TC004)
$ws.Range("B9").Value = "TC004"
$ws.Range("C9").Value = "MFP is powered on and in normal boot-up mode"
$ws.Range("D9").Value = "Verify SCC first registration URL can be configured"
$ws.Range("E9").Value = "1. Enter Self-diagnostic boot-up mode`n2. Set 08-3827 to a valid URL`n3. Reboot the MFP`n4. Verify the MFP attempts to connect to the specified URL"
$ws.Range("F9").Value = "MFP should attempt to connect to the specified URL for first registration"
$ws.Range("H9").Value = "Default value should be NULL string"

# Row 10 (TC005)
$ws.Range("B10").Value = "TC005"
$ws.Range("C10").Value = "MFP is powered on and in normal boot-up mode"
$ws.Range("D10").Value = "Verify URL forward setting can be configured"
$ws.Range("E10").Value = "1. Enter Self-diagnostic boot-up mode`n2. Set 08-3828 to 0 (Disable)`n3. Reboot the MFP`n4. Verify the MFP does not redirect to another URL"
$ws.Range("F10").Value = "MFP should not redirect to another URL when URL forward setting is disabled"
$ws.Range("H10").Value = ""

# Row 11 (TC006)
$ws.Range("B11").Value = "TC006"
$ws.Range("C11").Value = "SCC function is disabled"
$ws.Range("D11").Value = "Verify SCC Installation Report is printed when SCC function is enabled for the first time"
$ws.Range("E11").Value = "1. Enter Self-diagnostic boot-up mode`n2. Set 08-3820 from 0 (Disable) to 1 (Enable)`n3. Reboot the MFP`n4. Wait for MFP registration processing to complete"
$ws.Range("F11").Value = "Installation report should be automatically printed after MFP registration processing is executed"
$ws.Range("H11").Value = "Report should include: Model Name, Serial Number, IP Address, MAC Address, SCC Version, Date/Time"

# Row 12 (TC007)
$ws.Range("B12").Value = "TC007"
$ws.Range("C12").Value = "SCC function is already enabled"
$ws.Range("D12").Value = "Verify SCC Installation Report is not printed when SCC is already enabled"
$ws.Range("E12").Value = "1. With SCC already enabled, reboot the MFP`n2. Observe if installation report is printed"
$ws.Range("F12").Value = "Installation report should not be printed when SCC is already enabled"
$ws.Range("H12").Value = ""

# Row 13 (TC008)
$ws.Range("B13").Value = "TC008"
$ws.Range("C13").Value = "SCC function is enabled, network print restriction mode is set to `"only private`""
$ws.Range("D13").Value = "Verify SCC Installation Report behavior with print restrictions"
$ws.Range("E13").Value = "1. Set 08-9344 (network print restriction mode) to 1 (only private)`n2. Disable SCC function`n3. Enable SCC function`n4. Check if installation report is printed"
$ws.Range("F13").Value = "Installation report should not be printed, and job log should show error code 4221"
$ws.Range("H13").Value = ""

# Row 14 (TC009)
$ws.Range("B14").Value = "TC009"
$ws.Range("C14").Value = "SCC function is enabled"
$ws.Range("D14").Value = "Verify HTTPS communication protocol is used for server communication"
$ws.Range("E14").Value = "1. Enable network packet capture`n2. Trigger SCC communication with server`n3. Analyze captured packets"
$ws.Range("F14").Value = "All communication between MFP and server should use HTTPS protocol"
$ws.Range("H14").Value = ""

# Row 15 (TC010)
$ws.Range("B15").Value = "TC010"
$ws.Range("C15").Value = "SCC function is enabled, proxy server is configured"
$ws.Range("D15").Value = "Verify communication through proxy server"
$ws.Range("E15").Value = "1. Configure proxy server settings (08-3822 to 08-3826)`n2. Configure firewall to require proxy for outbound connections`n3. Trigger SCC communication`n4. Verify communication succeeds through proxy"
$ws.Range("F15").Value = "MFP should successfully communicate with SCC server through the configured proxy"
$ws.Range("H15").Value = ""

# Row 16 (TC011)
$ws.Range("B16").Value = "TC011"
$ws.Range("C16").Value = "SCC function is enabled, port 443 is blocked"
$ws.Range("D16").Value = "Verify automatic fallback to port 8443"
$ws.Range("E16").Value = "1. Block port 443 on the firewall`n2. Trigger SCC communication`n3. Monitor network traffic"
$ws.Range("F16").Value = "MFP should automatically use port 8443 when port 443 is not available"
$ws.Range("H16").Value = ""

# Row 17 (TC012)
$ws.Range("B17").Value = "TC012"
$ws.Range("C17").Value = "SCC function is enabled"
$ws.Range("D17").Value = "Verify SCC process does not start in special startup mode"
$ws.Range("E17").Value = "1. Boot the MFP in special startup mode`n2. Check if SCC process is running"
$ws.Range("F17").Value = "SCC process should not start in special startup mode"
$ws.Range("H17").Value = ""

# Row 18 (TC013)
$ws.Range("B18").Value = "TC013"
$ws.Range("C18").Value = "SCC function is enabled"
$ws.Range("D18").Value = "Verify regular communication loop operates according to schedule"
$ws.Range("E18").Value = "1. Set regular communication schedule to a specific time`n2. Monitor network traffic at the scheduled time"
$ws.Range("F18").Value = "MFP should initiate communication with the server at the scheduled time"
$ws.Range("H18").Value = "Default schedule is `"every day at 0:00`""

# Row 19 (TC014)
$ws.Range("B19").Value = "TC014"
$ws.Range("C19").Value = "SCC function is enabled"
$ws.Range("D19").Value = "Verify event notification loop operates when an event occurs"
$ws.Range("E19").Value = "1. Trigger a device error event`n2. Monitor network traffic"
$ws.Range("F19").Value = "MFP should immediately initiate communication with the server to notify the event"
$ws.Range("H19").Value = ""

# Row 20 (TC015)
$ws.Range("B20").Value = "TC015"
$ws.Range("C20").Value = "SCC function is enabled, MFP not registered"
$ws.Range("D20").Value = "Verify first-time MFP registration sequence"
$ws.Range("E20").Value = "1. Enable SCC on an unregistered MFP`n2. Monitor registration communication sequence"
$ws.Range("F20").Value = "MFP should send registration request, receive secret question, answer correctly, and complete registration"
$ws.Range("H20").Value = ""

# Row 21 (TC016)
$ws.Range("B21").Value = "TC016"
$ws.Range("C21").Value = "SCC function is enabled, MFP already registered"
$ws.Range("D21").Value = "Verify authentication sequence for registered MFP"
$ws.Range("E21").Value = "1. Trigger communication with server on a registered MFP`n2. Monitor authentication sequence"
$ws.Range("F21").Value = "MFP should connect using the token received during registration"
$ws.Range("H21").Value = ""

# Row 22 (TC017)
$ws.Range("B22").Value = "TC017"
$ws.Range("C22").Value = "SCC function is enabled, server is busy"
$ws.Range("D22").Value = "Verify retry behavior when server is busy"
$ws.Range("E22").Value = "1. Simulate server busy condition`n2. Trigger communication with server`n3. Monitor retry attempts"
$ws.Range("F22").Value = "MFP should retry connection with increasing intervals up to 3 times"
$ws.Range("H22").Value = ""

# Row 23 (TC018)
$ws.Range("B23").Value = "TC018"
$ws.Range("C23").Value = "SCC function is enabled"
$ws.Range("D23").Value = "Verify Check for Updates functionality"
$ws.Range("E23").Value = "1. Trigger regular communication cycle`n2. Monitor Check for Updates request`n3. Verify MFP sends correct parameters"
$ws.Range("F23").Value = "MFP should send correct parameters and process server response appropriately"
$ws.Range("H23").Value = ""

# Row 24 (TC019)
$ws.Range("B24").Value = "TC019"
$ws.Range("C24").Value = "SCC function is enabled, update package available"
$ws.Range("D24").Value = "Verify Download Package functionality"
$ws.Range("E24").Value = "1. Configure server to provide an update package`n2. Trigger regular communication cycle`n3. Monitor package download"
$ws.Range("F24").Value = "MFP should download the package and verify its hash value"
$ws.Range("H24").Value = ""

# Row 25 (TC020)
$ws.Range("B25").Value = "TC020"
$ws.Range("C25").Value = "SCC function is enabled, firmware update package downloaded"
$ws.Range("D25").Value = "Verify firmware update installation"
$ws.Range("E25").Value = "1. Download firmware update package`n2. Monitor installation process"
$ws.Range("F25").Value = "MFP should install the firmware update at the scheduled time and reboot"
$ws.Range("H25").Value = ""

# Row 26 (TC021)
$ws.Range("B26").Value = "TC021"
$ws.Range("C26").Value = "SCC function is enabled, policy violation data package downloaded"
$ws.Range("D26").Value = "Verify policy violation data installation"
$ws.Range("E26").Value = "1. Download policy violation data package`n2. Monitor installation process"
$ws.Range("F26").Value = "MFP should apply the policy settings correctly"
$ws.Range("H26").Value = ""

# Row 27 (TC022)
$ws.Range("B27").Value = "TC022"
$ws.Range("C27").Value = "SCC function is enabled, restore data package downloaded"
$ws.Range("D27").Value = "Verify restore data installation"
$ws.Range("E27").Value = "1. Download restore data package`n2. Monitor installation process"
$ws.Range("F27").Value = "MFP should restore the specified settings correctly"
$ws.Range("H27").Value = ""

# Row 28 (TC023)
$ws.Range("B28").Value = "TC023"
$ws.Range("C28").Value = "SCC function is enabled, clone file package downloaded"
$ws.Range("D28").Value = "Verify clone file installation"
$ws.Range("E28").Value = "1. Download clone file package`n2. Monitor installation process"
$ws.Range("F28").Value = "MFP should apply the clone settings correctly"
$ws.Range("H28").Value = ""

# Row 29 (TC024)
$ws.Range("B29").Value = "TC024"
$ws.Range("C29").Value = "SCC function is enabled, customized UI package downloaded"
$ws.Range("D29").Value = "Verify customized UI installation"
$ws.Range("E29").Value = "1. Download customized UI package`n2. Monitor installation process"
$ws.Range("F29").Value = "MFP should install the customized UI correctly"
$ws.Range("H29").Value = ""

# Row 30 (TC025)
$ws.Range("B30").Value = "TC025"
$ws.Range("C30").Value = "SCC function is enabled, application package downloaded"
$ws.Range("D30").Value = "Verify application installation"
$ws.Range("E30").Value = "1. Download application package`n2. Monitor installation process"
$ws.Range("F30").Value = "MFP should install the application correctly"
$ws.Range("H30").Value = ""

# Row 31 (TC026)
$ws.Range("B31").Value = "TC026"
$ws.Range("C31").Value = "SCC function is enabled, license package downloaded"
$ws.Range("D31").Value = "Verify license installation"
$ws.Range("E31").Value = "1. Download license package`n2. Monitor installation process"
$ws.Range("F31").Value = "MFP should install the license correctly"
$ws.Range("H31").Value = ""

# Row 32 (TC027)
$ws.Range("B32").Value = "TC027"
$ws.Range("C32").Value = "SCC function is enabled, custom user paper type package downloaded"
$ws.Range("D32").Value = "Verify custom user paper type installation"
$ws.Range("E32").Value = "1. Download custom user paper type package`n2. Monitor installation process"
$ws.Range("F32").Value = "MFP should install the custom user paper type correctly"
$ws.Range("H32").Value = ""

# Row 33 (TC028)
$ws.Range("B33").Value = "TC028"
$ws.Range("C33").Value = "SCC function is enabled, update package downloaded"
$ws.Range("D33").Value = "Verify Update Status notification"
$ws.Range("E33").Value = "1. Install an update package`n2. Monitor Update Status notification"
$ws.Range("F33").Value = "MFP should send correct update status to server"
$ws.Range("H33").Value = ""

# Row 34 (TC029)
$ws.Range("B34").Value = "TC029"
$ws.Range("C34").Value = "SCC function is enabled"
$ws.Range("D34").Value = "Verify Send Baseline Data functionality"
$ws.Range("E34").Value = "1. Trigger regular communication cycle`n2. Monitor Baseline Data transmission"
$ws.Range("F34").Value = "MFP should collect and send all required baseline data to server"
$ws.Range("H34").Value = ""

# Row 35 (TC030)
$ws.Range("B35").Value = "TC030"
$ws.Range("C35").Value = "SCC function is enabled"
$ws.Range("D35").Value = "Verify Send Regular Data functionality"
$ws.Range("E35").Value = "1. Trigger regular communication cycle`n2. Monitor Regular Data transmission"
$ws.Range("F35").Value = "MFP should collect and send all required regular data to server"
$ws.Range("H35").Value = ""

# Row 36 (TC031)
$ws.Range("B36").Value = "TC031"
$ws.Range("C36").Value = "SCC function is enabled, device error occurs"
$ws.Range("D36").Value = "Verify Send Device Error functionality"
$ws.Range("E36").Value = "1. Trigger a device error`n2. Monitor Device Error notification"
$ws.Range("F36").Value = "MFP should send error details to server immediately"
$ws.Range("H36").Value = ""

# Row 37 (TC032)
$ws.Range("B37").Value = "TC032"
$ws.Range("C37").Value = "SCC function is enabled, same device error occurs multiple times"
$ws.Range("D37").Value = "Verify duplicate error handling"
$ws.Range("E37").Value = "1. Trigger the same device error multiple times`n2. Monitor Device Error notifications"
$ws.Range("F37").Value = "MFP should not send duplicate error notifications for the same error"
$ws.Range("H37").Value = ""

# Row 38 (TC033)
$ws.Range("B38").Value = "TC033"
$ws.Range("C38").Value = "SCC function is enabled, device error is resolved"
$ws.Range("D38").Value = "Verify error resolution notification"
$ws.Range("E38").Value = "1. Trigger a device error`n2. Resolve the error without power off`n3. Monitor Device Error notifications"
$ws.Range("F38").Value = "MFP should notify server of error resolution with `"-`" prefix (e.g., -D102)"
$ws.Range("H38").Value = ""

# Row 39 (TC034)
$ws.Range("B39").Value = "TC034"
$ws.Range("C39").Value = "SCC function is enabled, server requests service files"
$ws.Range("D39").Value = "Verify Send Service File functionality"
$ws.Range("E39").Value = "1. Configure server to request service files`n2. Trigger device error notification`n3. Monitor Service File transmission"
$ws.Range("F39").Value = "MFP should send ZIP file containing service files to server"
$ws.Range("H39").Value = ""

# Row 40 (TC035)
$ws.Range("B40").Value = "TC035"
$ws.Range("C40").Value = "SCC function is enabled, MFP in Super Sleep"
$ws.Range("D40").Value = "Verify MFP wakes from Super Sleep for scheduled communication"
$ws.Range("E40").Value = "1. Put MFP in Super Sleep state`n2. Wait for scheduled communication time`n3. Observe MFP behavior"
$ws.Range("F40").Value = "MFP should wake from Super Sleep, perform communication, then return to Super Sleep"
$ws.Range("H40").Value = ""

# Row 41 (TC036)
$ws.Range("B41").Value = "TC036"
$ws.Range("C41").Value = "SCC function is enabled, MFP about to enter Super Sleep"
$ws.Range("D41").Value = "Verify MFP does not enter Super Sleep during communication"
$ws.Range("E41").Value = "1. Trigger SCC communication`n2. Attempt to put MFP in Super Sleep state`n3. Observe MFP behavior"
$ws.Range("F41").Value = "MFP should not enter Super Sleep during communication cycle"
$ws.Range("H41").Value = ""

# Row 42 (TC037)
$ws.Range("B42").Value = "TC037"
$ws.Range("C42").Value = "SCC function is enabled, persistent policy configured"
$ws.Range("D42").Value = "Verify persistent policy check functionality"
$ws.Range("E42").Value = "1. Configure persistent policy with interval`n2. Change settings locally`n3. Wait for persistent policy check`n4. Verify settings"
$ws.Range("F42").Value = "Settings should be restored to policy values after local changes"
$ws.Range("H42").Value = ""

# Row 43 (TC038)
$ws.Range("B43").Value = "TC038"
$ws.Range("C43").Value = "SCC function is enabled, time-based device state configured"
$ws.Range("D43").Value = "Verify time-based device state functionality"
$ws.Range("E43").Value = "1. Configure time-based values for settings`n2. Observe setting changes at specified times"
$ws.Range("F43").Value = "Settings should change automatically at the specified times"
$ws.Range("H43").Value = ""

# Row 44 (TC039)
$ws.Range("B44").Value = "TC039"
$ws.Range("C44").Value = "SCC function is enabled, first-time connection"
$ws.Range("D44").Value = "Verify IP Redirect functionality"
$ws.Range("E44").Value = "1. Configure MFP for first-time connection`n2. Monitor GetRedirectURL request`n3. Verify MFP follows redirect"
$ws.Range("F44").Value = "MFP should request redirect URL and connect to appropriate regional server"
$ws.Range("H44").Value = ""

# Row 45 (TC040)
$ws.Range("B45").Value = "TC040"
$ws.Range("C45").Value = "SCC function is enabled, communication error occurs"
$ws.Range("D45").Value = "Verify retry behavior for communication errors"
$ws.Range("E45").Value = "1. Simulate network error during communication`n2. Monitor retry attempts`n3. Verify retry mode behavior"
$ws.Range("F45").Value = "MFP should retry after 60 seconds, then enter retry mode with daily attempts"
$ws.Range("H45").Value = ""

# Row 46 (TC041)
$ws.Range("B46").Value = "TC041"
$ws.Range("C46").Value = "SCC function is enabled, HTTP error occurs"
$ws.Range("D46").Value = "Verify retry behavior for HTTP errors"
$ws.Range("E46").Value = "1. Simulate HTTP error from server`n2. Monitor retry attempts`n3. Verify retry mode behavior"
$ws.Range("F46").Value = "MFP should retry after 60 seconds, then enter retry mode with daily attempts"
$ws.Range("H46").Value = ""

# Row 47 (TC042)
$ws.Range("B47").Value = "TC042"
$ws.Range("C47").Value = "SCC function is enabled, power failure during communication"
$ws.Range("D47").Value = "Verify recovery from power failure"
$ws.Range("E47").Value = "1. Trigger SCC communication`n2. Simulate power failure during communication`n3. Restore power`n4. Observe MFP behavior"
$ws.Range("F47").Value = "MFP should start normal communication cycle after power restoration"
$ws.Range("H47").Value = ""
$ws.Range("G47").Value = ""

# Row 48 (TC043)
$ws.Range("B48").Value = "TC043"
$ws.Range("C48").Value = "SCC function is enabled, HDD full during package download"
$ws.Range("D48").Value = "Verify behavior when HDD is full"
$ws.Range("E48").Value = "1. Fill MFP HDD to capacity`n2. Trigger package download`n3. Observe MFP behavior"
$ws.Range("F48").Value = "MFP should delete downloaded data and exit communication cycle"
$ws.Range("H48").Value = ""
$ws.Range("G48").Value = ""

# Row 49 (TC044)
$ws.Range("B49").Value = "TC044"
$ws.Range("C49").Value = "SCC function is enabled, hash value mismatch"
$ws.Range("D49").Value = "Verify behavior when package hash value is incorrect"
$ws.Range("E49").Value = "1. Simulate package with incorrect hash value`n2. Trigger package download`n3. Observe MFP behavior"
$ws.Range("F49").Value = "MFP should delete downloaded data and exit communication cycle"
$ws.Range("H49").Value = ""
$ws.Range("G49").Value = ""

# Row 50 (TC045)
$ws.Range("B50").Value = "TC045"
$ws.Range("C50").Value = "SCC function is enabled, unzip failure"
$ws.Range("D50").Value = "Verify behavior when package unzip fails"
$ws.Range("E50").Value = "1. Simulate package that cannot be unzipped`n2. Trigger package download`n3. Observe MFP behavior"
$ws.Range("F50").Value = "MFP should delete unzipped data and exit communication cycle"
$ws.Range("H50").Value = ""
$ws.Range("G50").Value = ""

# Row 51 (TC046)
$ws.Range("B51").Value = "TC046"
$ws.Range("C51").Value = "SCC function is enabled, installation failure"
$ws.Range("D51").Value = "Verify behavior when installation fails"
$ws.Range("E51").Value = "1. Simulate installation failure`n2. Observe MFP behavior"
$ws.Range("F51").Value = "MFP should notify update result to server and schedule regular communication after 1 hour"
$ws.Range("H51").Value = ""
$ws.Range("G51").Value = ""

# Row 52 (TC047)
$ws.Range("B52").Value = "TC047"
$ws.Range("C52").Value = "SCC function is enabled, SSD model"
$ws.Range("D52").Value = "Verify firmware update behavior on SSD model"
$ws.Range("E52").Value = "1. Attempt firmware update on SSD model`n2. Observe MFP behavior"
$ws.Range("F52").Value = "Firmware update should fail, and MFP should notify server of failure"
$ws.Range("H52").Value = ""
$ws.Range("G52").Value = ""

# Row 53 (TC048)
$ws.Range("B53").Value = "TC048"
$ws.Range("C53").Value = "SCC function is enabled, multiple updates with partial failure"
$ws.Range("D53").Value = "Verify behavior when some updates fail"
$ws.Range("E53").Value = "1. Configure multiple updates with one set to fail`n2. Trigger update installation`n3. Observe MFP behavior"
$ws.Range("F53").Value = "MFP should continue processing all updates and report results to server"
$ws.Range("H53").Value = ""
$ws.Range("G53").Value = ""

# Row 54 (TC049)
$ws.Range("B54").Value = "TC049"
$ws.Range("C54").Value = "SCC function is enabled"
$ws.Range("D54").Value = "Verify panel message display during SCC processing"
$ws.Range("E54").Value = "1. Trigger SCC communication`n2. Observe panel display"
$ws.Range("F54").Value = "Panel should display `"Service in progress. Please do not turn off: XX`" with appropriate status code"
$ws.Range("H54").Value = ""
$ws.Range("G54").Value = ""

# Row 55 (TC050)
$ws.Range("B55").Value = "TC050"
$ws.Range("C55").Value = "SCC function is enabled, firmware update in progress"
$ws.Range("D55").Value = "Verify popup window display during firmware update"
$ws.Range("E55").Value = "1. Trigger firmware update`n2. Observe panel display"
$ws.Range("F55").Value = "Panel should display popup window with sand clock and progress bar"
$ws.Range("H55").Value = ""
$ws.Range("G55").Value = ""

# Row 56 (TC051)
$ws.Range("B56").Value = "TC051"
$ws.Range("C56").Value = "SCC function is enabled, RDMS (BBR2) function is enabled"
$ws.Range("D56").Value = "Verify SCC and RDMS exclusivity"
$ws.Range("E56").Value = "1. Enable SCC function`n2. Enable RDMS (BBR2) function`n3. Attempt to use RDMS"
$ws.Range("F56").Value = "RDMS should not work when SCC function is enabled"
$ws.Range("H56").Value = "This limitation is removed from L6.02/L6.03"
$ws.Range("G56").Value = ""

# Row 57 (TC052)
$ws.Range("B57").Value = "TC052"
$ws.Range("C57").Value = "SCC function is enabled, service UI in use"
$ws.Range("D57").Value = "Verify exclusivity between SCC and service UI"
$ws.Range("E57").Value = "1. Access service UI`n2. Trigger SCC communication`n3. Attempt to use service UI"
$ws.Range("F57").Value = "Service UI should not be usable while device is communicating with SCC server"
$ws.Range("H57").Value = ""
$ws.Range("G57").Value = ""

# Row 58 (TC053)
$ws.Range("B58").Value = "TC053"
$ws.Range("C58").Value = "SCC function is enabled, invalid URL in 08-3827"
$ws.Range("D58").Value = "Verify behavior with invalid URL"
$ws.Range("E58").Value = "1. Set 08-3827 to an invalid URL format`n2. Trigger SCC communication`n3. Observe MFP behavior"
$ws.Range("F58").Value = "Network error should occur, and MFP should not access NA server"
$ws.Range("H58").Value = ""
$ws.Range("G58").Value = ""

# Row 59 (TC054)
$ws.Range("B59").Value = "TC054"
$ws.Range("C59").Value = "SCC function is enabled, URL without http/https prefix"
$ws.Range("D59").Value = "Verify behavior with URL missing protocol"
$ws.Range("E59").Value = "1. Set 08-3827 to URL without http/https prefix`n2. Trigger SCC communication`n3. Observe MFP behavior"
$ws.Range("F59").Value = "Network error should occur, and MFP should not access NA server"
$ws.Range("H59").Value = ""
$ws.Range("G59").Value = ""

# Row 60 (TC055)
$ws.Range("B60").Value = "TC055"
$ws.Range("C60").Value = "SCC function is enabled, HTTPS URL without imported certificate"
$ws.Range("D60").Value = "Verify behavior with HTTPS URL lacking certificate"
$ws.Range("E60").Value = "1. Set 08-3827 to HTTPS URL without importing certificate`n2. Trigger SCC communication`n3. Observe MFP behavior"
$ws.Range("F60").Value = "Communication should fail due to missing certificate"
$ws.Range("H60").Value = ""
$ws.Range("G60").Value = ""

# Row 61 (TC056)
$ws.Range("B61").Value = "TC056"
$ws.Range("C61").Value = "SCC function is enabled, power failure during firmware update wait"
$ws.Range("D61").Value = "Verify firmware update behavior after power failure"
$ws.Range("E61").Value = "1. Schedule firmware update`n2. Power off MFP before scheduled time`n3. Power on MFP after scheduled time`n4. Observe MFP behavior"
$ws.Range("F61").Value = "For L6.15+: Firmware update should be rescheduled for same time next day`nBefore L6.15: Firmware update should execute 1 hour after reboot"
$ws.Range("H61").Value = ""
$ws.Range("G61").Value = ""

# Row 62 (TC057)
$ws.Range("B62").Value = "TC057"
$ws.Range("C62").Value = "SCC function is enabled, hibernation about to occur"
$ws.Range("D62").Value = "Verify MFP does not enter hibernation during communication"
$ws.Range("E62").Value = "1. Trigger SCC communication`n2. Attempt to put MFP in hibernation state`n3. Observe MFP behavior"
$ws.Range("F62").Value = "MFP should not enter hibernation during communication cycle"
$ws.Range("H62").Value = ""
$ws.Range("G62").Value = ""

# Row 63 (TC058)
$ws.Range("B63").Value = "TC058"
$ws.Range("C63").Value = "SCC function is enabled, persistent policy expiration configured"
$ws.Range("D63").Value = "Verify persistent policy expiration functionality"
$ws.Range("E63").Value = "1. Configure persistent policy with expiration`n2. Wait for expiration period`n3. Change settings locally`n4. Verify if settings are restored"
$ws.Range("F63").Value = "Settings should not be restored after persistent policy expiration"
$ws.Range("H63").Value = ""
$ws.Range("G63").Value = ""

# Row 64 (TC059)
$ws.Range("B64").Value = "TC059"
$ws.Range("C64").Value = "SCC function is enabled, event occurs and cancels immediately"
$ws.Range("D64").Value = "Verify event notification behavior for quickly canceled events"
$ws.Range("E64").Value = "1. Trigger an event that cancels immediately`n2. Observe if event is notified to server"
$ws.Range("F64").Value = "Event may not be notified if SCC error notification process is in busy state"
$ws.Range("H64").Value = ""
$ws.Range("G64").Value = ""

# Row 65 (TC060)
$ws.Range("B65").Value = "TC060"
$ws.Range("C65").Value = "SCC function is enabled, F-category error occurs"
$ws.Range("D65").Value = "Verify behavior when F-category error occurs"
$ws.Range("E65").Value = "1. Trigger F-category error`n2. Observe if error is sent to server"
$ws.Range("F65").Value = "Device should not be able to send error to server due to network unavailability"
$ws.Range("H65").Value = ""
$ws.Range("G65").Value = ""

